$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (month strings). Assigning a "Month Year" looking
# string straight to .Value makes Excel auto-convert it into a date serial
# number, so instead enter it as a text-formula ("=""April 2024""") and then
# convert the cell back down to a plain text value via copy/paste-values -
# this avoids touching NumberFormat (which would otherwise mint a new,
# unwanted cell style).
$ws.Range("A1").Formula = '="April 2024"'
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4163) | Out-Null

$ws.Range("G1").Formula = '="May 2024"'
$ws.Range("G1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0

# Update numeric row 2 values
$ws.Range("A2").Value = 1.602
$ws.Range("B2").Value = 0.278
$ws.Range("C2").Value = -0.092
$ws.Range("D2").Value = -0.061
$ws.Range("E2").Value = -0.026
$ws.Range("F2").Value = 0.232
$ws.Range("G2").Value = 1.934
